$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.210.02'

$ws.Cells.Item(3, 4).Value = '1.969.47'
$ws.Cells.Item(3, 5).Value = '  -6.35%  '

$ws.Cells.Item(4, 4).Value = "'1.010"
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  +0.65%  '

$ws.Cells.Item(5, 4).Value = "'329.32"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -3.95%  '

$ws.Cells.Item(6, 4).Value = "'1.010"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +0.77%  '

$ws.Cells.Item(7, 4).Value = "'0.4977"
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -5.73%  '

$ws.Cells.Item(8, 4).Value = "'0.4240"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -3.99%  '

$ws.Cells.Item(9, 4).Value = "'54.45"
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -1.06%  '

$ws.Cells.Item(10, 4).Value = "'0.09190"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -1.81%  '

$ws.Cells.Item(11, 4).Value = "'1.102"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -6.25%  '

$ws.Cells.Item(12, 4).Value = "'23.25"
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  -6.20%  '

$ws.Cells.Item(13, 4).Value = "'7.950"
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  -6.94%  '

$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(14, 4).Value = "'6.473"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -6.02%  '

$ws.Cells.Item(15, 2).Value = 'WrappedEther'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(15, 4).Value = '1.934.13'
$ws.Cells.Item(15, 5).Value = '  -11.94%  '

$ws.Cells.Item(16, 4).Value = "'1.012"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +0.81%  '

$ws.Cells.Item(17, 4).Value = "'0.00001109"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -4.29%  '

$ws.Cells.Item(18, 4).Value = "'91.87"
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  -9.48%  '

$ws.Cells.Item(19, 4).Value = "'0.06698"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -0.22%  '

$ws.Cells.Item(20, 4).Value = "'19.35"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -8.76%  '

$ws.Cells.Item(21, 4).Value = "'1.008"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +0.48%  '

$ws.Cells.Item(22, 4).Value = "'5.944"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -6.82%  '

$ws.Cells.Item(23, 4).Value = '29.206.46'
$ws.Cells.Item(23, 5).Value = '  -3.94%  '

$ws.Cells.Item(24, 4).Value = "'12.03"
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -3.61%  '

$ws.Cells.Item(25, 4).Value = "'2.297"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -0.93%  '

$ws.Cells.Item(26, 4).Value = "'20.74"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -5.10%  '

$ws.Cells.Item(27, 4).Value = "'157.10"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -3.20%  '

$ws.Cells.Item(28, 4).Value = "'6.247"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -12.17%  '

$ws.Cells.Item(29, 4).Value = "'2.283"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -9.37%  '

$ws.Cells.Item(30, 4).Value = "'127.26"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -4.93%  '

$ws.Cells.Item(31, 4).Value = "'1.048"
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -8.65%  '

$ws.Cells.Item(32, 4).Value = "'0.09879"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -6.39%  '

$ws.Cells.Item(33, 4).Value = "'1.552"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -7.98%  '

$ws.Cells.Item(34, 4).Value = "'5.827"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -7.45%  '

$ws.Cells.Item(35, 4).Value = "'3.680"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -5.19%  '

$ws.Cells.Item(36, 4).Value = "'0.02442"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -7.41%  '

$ws.Cells.Item(37, 4).Value = "'9.082"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -11.11%  '

$ws.Cells.Item(38, 4).Value = "'0.06373"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -5.95%  '

$ws.Cells.Item(39, 4).Value = "'1.293"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -3.83%  '

$ws.Cells.Item(40, 4).Value = "'0.6500"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -7.25%  '

$ws.Cells.Item(41, 4).Value = "'11.52"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -9.05%  '

$ws.Cells.Item(42, 4).Value = "'0.2013"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -9.37%  '

$ws.Cells.Item(43, 4).Value = "'1.009"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +0.72%  '

$ws.Cells.Item(44, 4).Value = "'0.6262"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -7.89%  '

$ws.Cells.Item(45, 4).Value = "'13.62"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -6.08%  '

$ws.Cells.Item(46, 4).Value = "'2.192"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -5.90%  '

$ws.Cells.Item(47, 4).Value = "'1.279"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -4.62%  '

$ws.Cells.Item(48, 4).Value = "'3.492"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -3.85%  '

$ws.Cells.Item(49, 4).Value = "'0.00000000339"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -1.68%  '

$ws.Cells.Item(50, 4).Value = "'0.06889"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -5.95%  '

$ws.Cells.Item(51, 4).Value = "'1.116"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -8.30%  '
